$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns to fill with "x" across rows 3, 5 and 19 (matching the TODO checklist
# being marked done for the extra columns, except D and O which are untouched).
$cols = @("C", "E", "F", "H", "I", "J", "K", "L", "M", "N")
$rows = @(3, 5, 19)

foreach ($row in $rows) {
    foreach ($col in $cols) {
        $ws.Range("$col$row").Value = "x"
    }
}

# Move the active selection to P5 as recorded in the saved view state.
$ws.Range("P5").Select()
